# Insert a new data row at row 138 (pushing the existing rows 138-200 down
# to 139-201) for the "Vega Modelo de Temuco - Espinaca" weekly price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 138:200 down by one row.
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new weekly record.
$ws.Cells.Item(138, 1).Value  = 10
$ws.Cells.Item(138, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(138, 3).Value  = "La Araucanía"
$ws.Cells.Item(138, 4).Value  = 44839
$ws.Cells.Item(138, 5).Value  = 9
$ws.Cells.Item(138, 6).Value  = 100112012
$ws.Cells.Item(138, 7).Value  = "Espinaca"
$ws.Cells.Item(138, 8).Value  = "Sin especificar"
$ws.Cells.Item(138, 9).Value  = "Primera"
$ws.Cells.Item(138, 10).Value = 50
$ws.Cells.Item(138, 11).Value = 10000
$ws.Cells.Item(138, 12).Value = 10000
$ws.Cells.Item(138, 13).Value = 10000
$ws.Cells.Item(138, 14).Value = "$/docena de atados"
$ws.Cells.Item(138, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(138, 16).Value = 3333
$ws.Cells.Item(138, 17).Value = 3
$ws.Cells.Item(138, 18).Value = "Hortaliza"
